$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell "Save" in H1, reusing the same formatting as the
# existing header row (copy format from G1 so it shares the style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for rows 2-9
$values = @(1, 1, 0, 0, 1, 0, 0, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
